$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue 'D2' '298.19'
Set-TextValue 'E2' '-3.46%'
Set-TextValue 'D3' '31.93'
Set-TextValue 'E3' '-0.81%'
Set-TextValue 'D4' '5.089'
Set-TextValue 'E4' '-4.75%'
Set-TextValue 'D5' '0.07531'
Set-TextValue 'E5' '0.68%'
Set-TextValue 'D6' '7.750'
Set-TextValue 'E6' '-0.52%'
Set-TextValue 'D7' '1.725'
Set-TextValue 'E7' '9.19%'
Set-TextValue 'D8' '3.794'
Set-TextValue 'E8' '3.31%'
Set-TextValue 'D9' '0.9290'
Set-TextValue 'E9' '2.05%'
Set-TextValue 'D10' '0.1702'
Set-TextValue 'E10' '1.57%'
Set-TextValue 'D11' '0.07457'
Set-TextValue 'E11' '-1.46%'
Set-TextValue 'D12' '0.07927'
Set-TextValue 'E12' '-1.92%'
Set-TextValue 'D13' '0.03057'
Set-TextValue 'E13' '1.15%'
Set-TextValue 'D14' '0.09892'
Set-TextValue 'E14' '0.37%'
Set-TextValue 'D15' '0.001487'
Set-TextValue 'E15' '-2.65%'
Set-TextValue 'D16' '0.006455'
Set-TextValue 'E16' '-1.35%'
Set-TextValue 'D17' '3.465'
Set-TextValue 'E17' '-0.73%'
Set-TextValue 'D18' '2.222'
Set-TextValue 'E18' '-0.67%'
Set-TextValue 'E19' '0.44%'
Set-TextValue 'D20' '0.1326'
Set-TextValue 'E20' '-0.57%'
Set-TextValue 'D21' '4.556'
Set-TextValue 'E21' '9.03%'
Set-TextValue 'D22' '0.04657'
Set-TextValue 'E22' '2.67%'
Set-TextValue 'D23' '0.1558'
Set-TextValue 'E23' '-3.69%'
Set-TextValue 'E24' '0.50%'
Set-TextValue 'E25' '-1.99%'
Set-TextValue 'E26' '7.80%'
Set-TextValue 'E27' '6.94%'
Set-TextValue 'D39' '0.01678'
Set-TextValue 'E39' '-2.22%'
Set-TextValue 'D40' '0.04540'
Set-TextValue 'E40' '0.20%'
Set-TextValue 'D41' '0.007042'
Set-TextValue 'E41' '-1.91%'
Set-TextValue 'D42' '0.1327'
Set-TextValue 'E42' '-2.55%'
Set-TextValue 'D43' '0.002059'
Set-TextValue 'E43' '-8.76%'
Set-TextValue 'D44' '0.01172'
Set-TextValue 'E44' '-14.31%'
Set-TextValue 'D45' '0.00005984'
Set-TextValue 'E45' '-3.45%'
Set-TextValue 'D46' '1.918'
Set-TextValue 'E46' '1.35%'
Set-TextValue 'E47' '-0.05%'
